$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-09-16 Tuesday" "2025-09-17 Wednesday"

Replace-Text "605×4=" "854×7="
Replace-Text "573×9=" "372×6="
Replace-Text "241×9=" "342×4="
Replace-Text "729×5=" "685×5="
Replace-Text "953×7=" "595×4="

Replace-Text "364×7=" "638×4="
Replace-Text "217×8=" "101×4="
Replace-Text "578×2=" "238×2="
Replace-Text "190×6=" "915×7="
Replace-Text "952×9=" "143×8="

Replace-Text "460×3=" "626×9="
Replace-Text "786×7=" "542×4="
Replace-Text "163×8=" "433×8="
Replace-Text "556×4=" "373×7="
Replace-Text "214×4=" "740×8="

Replace-Text "353×6=" "157×9="
Replace-Text "769×6=" "411×3="
Replace-Text "365×8=" "208×2="
Replace-Text "309×6=" "708×6="
Replace-Text "660×5=" "355×6="

Replace-Text "890×6=" "849×7="
Replace-Text "554×4=" "465×3="
Replace-Text "379×7=" "457×6="
Replace-Text "498×2=" "684×9="
Replace-Text "112×3=" "691×3="

Write-Output "Done"
